$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values would otherwise be
# auto-coerced to numbers by Excel (so they stay text, matching the source).
$textCells = @("D5", "D6", "D7", "D12", "D13", "D16", "D17", "D19", "D20", "D21", "D22", "D24", "D28", "D30", "D32", "D33", "D35", "D36", "D38", "D40", "D41", "D42", "D44", "D45", "D47", "D48", "D49")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values
$ws.Range("D2").Value = "58.238.80"
$ws.Range("E2").Value = "  +0.10%  "
$ws.Range("D3").Value = "2.524.03"
$ws.Range("E3").Value = "  +1.82%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "522.11"
$ws.Range("E5").Value = "  +0.55%  "
$ws.Range("D6").Value = "133.02"
$ws.Range("E6").Value = "  -1.14%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("E8").Value = "  +1.19%  "
$ws.Range("D9").Value = "2.522.81"
$ws.Range("E9").Value = "  +1.25%  "
$ws.Range("E10").Value = "  -0.99%  "
$ws.Range("E11").Value = "  -1.51%  "
$ws.Range("D12").Value = "5.16"
$ws.Range("E12").Value = "  -2.70%  "
$ws.Range("D13").Value = "0.332"
$ws.Range("E13").Value = "  -2.42%  "
$ws.Range("D14").Value = "2.970.05"
$ws.Range("E14").Value = "  +1.78%  "
$ws.Range("D15").Value = "58.266.91"
$ws.Range("E15").Value = "  +0.30%  "
$ws.Range("D16").Value = "22.11"
$ws.Range("E16").Value = "  -0.07%  "
$ws.Range("D17").Value = "0.0000135"
$ws.Range("E17").Value = "  -0.38%  "
$ws.Range("D18").Value = "2.523.28"
$ws.Range("E18").Value = "  +1.53%  "
$ws.Range("D19").Value = "10.67"
$ws.Range("E19").Value = "  +0.06%  "
$ws.Range("D20").Value = "322.16"
$ws.Range("E20").Value = "  +0.33%  "
$ws.Range("D21").Value = "4.17"
$ws.Range("E21").Value = "  -0.53%  "
$ws.Range("D22").Value = "6.14"
$ws.Range("E22").Value = "  +6.89%  "
$ws.Range("E23").Value = "  +0.22%  "
$ws.Range("D24").Value = "64.50"
$ws.Range("E24").Value = "  +0.44%  "
$ws.Range("E25").Value = "  -1.11%  "
$ws.Range("E26").Value = "  +0.26%  "
$ws.Range("E27").Value = "  -1.03%  "
$ws.Range("D28").Value = "7.40"
$ws.Range("E28").Value = "  +0.11%  "
$ws.Range("E29").Value = "  +0.34%  "
$ws.Range("D30").Value = "168.76"
$ws.Range("E30").Value = "  -0.71%  "
$ws.Range("E31").Value = "  +1.40%  "
$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").Value = "1.19"
$ws.Range("E32").Value = "  +0.18%  "
$ws.Range("B33").Value = "Aptos"
$ws.Range("C33").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D33").Value = "6.30"
$ws.Range("E33").Value = "  -1.05%  "
$ws.Range("D35").Value = "0.997"
$ws.Range("E35").Value = "  +0.13%  "
$ws.Range("D36").Value = "18.18"
$ws.Range("E36").Value = "  +0.31%  "
$ws.Range("E37").Value = "  -6.41%  "
$ws.Range("D38").Value = "3.92"
$ws.Range("E38").Value = "  -2.77%  "
$ws.Range("E39").Value = "  +1.04%  "
$ws.Range("D40").Value = "36.47"
$ws.Range("E40").Value = "  -0.31%  "
$ws.Range("D41").Value = "0.772"
$ws.Range("E41").Value = "  -3.68%  "
$ws.Range("D42").Value = "276.34"
$ws.Range("E42").Value = "  +0.05%  "
$ws.Range("E43").Value = "  +0.01%  "
$ws.Range("D44").Value = "130.23"
$ws.Range("E44").Value = "  +5.29%  "
$ws.Range("D45").Value = "4.99"
$ws.Range("E45").Value = "  -2.98%  "
$ws.Range("E46").Value = "  +0.15%  "
$ws.Range("D47").Value = "0.0917"
$ws.Range("E47").Value = "  +0.78%  "
$ws.Range("D48").Value = "0.0500"
$ws.Range("E48").Value = "  +1.72%  "
$ws.Range("D49").Value = "17.70"
$ws.Range("E50").Value = "  +0.15%  "
$ws.Range("E51").Value = "  -0.95%  "

# Restore default (Normal) style on the cells we temporarily reformatted,
# now that the text value has "stuck" -- keeps cell styling identical to
# the original workbook.
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
